$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.487.50"
$ws.Range("E2").Value = "  +3.63%  "

$ws.Range("D3").Value = "2.427.74"
$ws.Range("E3").Value = "  +2.81%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.69"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.69"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +6.52%  "

$ws.Range("E7").Value = "  +2.09%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  +5.62%  "

$ws.Range("E10").Value = "  +4.03%  "

$ws.Range("E11").Value = "  +1.99%  "

$ws.Range("E12").Value = "  +1.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.83"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.24%  "

$ws.Range("E14").Value = "  +3.47%  "

$ws.Range("D15").Value = "2.806.49"
$ws.Range("E15").Value = "  +2.83%  "

$ws.Range("D16").Value = "2.420.64"
$ws.Range("E16").Value = "  +2.41%  "

$ws.Range("E17").Value = "  +5.40%  "

$ws.Range("D18").Value = "44.426.98"
$ws.Range("E18").Value = "  +3.64%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.43"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.27%  "

$ws.Range("E20").Value = "  +2.45%  "

$ws.Range("E21").Value = "  +2.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.94"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.64%  "

$ws.Range("E23").Value = "  +2.61%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.28"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +4.12%  "

$ws.Range("E25").Value = "  +2.28%  "

$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.20"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.26%  "

$ws.Range("E28").Value = "  -4.24%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.63"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.25%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.31"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +5.71%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "48.49"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.26%  "

$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.122"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +17.01%  "

$ws.Range("B33").Value = "Celestia"
$ws.Range("C33").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.52"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +12.57%  "

$ws.Range("E34").Value = "  +3.53%  "

$ws.Range("E35").Value = "  +0.27%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0767"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +6.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.54"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +4.56%  "

$ws.Range("E38").Value = "  +3.33%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.90"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.70%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "126.59"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.79%  "

$ws.Range("E41").Value = "  +1.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.77"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.65%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.16"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -6.36%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0289"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.65%  "

$ws.Range("D45").Value = "1.950.44"
$ws.Range("E45").Value = "  +0.95%  "

$ws.Range("E46").Value = "  +1.75%  "

$ws.Range("E47").Value = "  +8.85%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.62"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +4.66%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.67"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +10.41%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.37"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.46%  "

$ws.Range("E51").Value = "  +2.77%  "
